$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 9 data columns (A..I) get logically permuted: each old column's whole
# A1:x5 block moves to a new column position (single 9-cycle permutation),
# and the row-1 slug label is replaced by a human-readable (accented) title.
# Use a staging area (columns K..S) so every source column is copied before
# any destination is overwritten.

$oldCols    = @("A","B","C","D","E","F","G","H","I")
$stageCols  = @("K","L","M","N","O","P","Q","R","S")
$newColFor  = @{ "A"="C"; "B"="G"; "C"="F"; "D"="B"; "E"="D"; "F"="I"; "G"="H"; "H"="A"; "I"="E" }

# 1) Stage a copy of every old column (rows 1-5) so none of the upcoming
#    writes into A..I can clobber a still-to-be-read source column.
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $src = $ws.Range($oldCols[$i] + "1:" + $oldCols[$i] + "5")
    $dst = $ws.Range($stageCols[$i] + "1:" + $stageCols[$i] + "5")
    $src.Copy($dst)
}

# 2) Clear the live A..I area completely (values + styles) before laying the
#    reordered columns back down, since a couple of target cells (E5, H5)
#    must end up with no cell at all.
$ws.Range("A1:I5").ClearContents()
$ws.Range("A1:I5").ClearFormats()

# 3) Copy each staged column into its new destination column, preserving
#    the original cell style (s="1") that Copy()/paste carries along.
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $old = $oldCols[$i]
    $new = $newColFor[$old]
    $src = $ws.Range($stageCols[$i] + "1:" + $stageCols[$i] + "5")
    $dst = $ws.Range($new + "1:" + $new + "5")
    $src.Copy($dst)
}

# 4) The mapping-file cells for the measure column (now E) and the
#    territory/URI column (now H) don't exist in row 5 any more.
$ws.Range("E5").ClearContents()
$ws.Range("H5").ClearContents()

# 5) Drop the staging columns entirely.
$ws.Range("K1:S5").Clear()

# 6) Row 1 gets the new human-readable (and accented) labels, replacing the
#    old machine slugs that were carried over by the column move above.
$ws.Range("A1").Value = "Poca limpieza"
$ws.Range("B1").Value = "Pocas zonas verdes"
$ws.Range("C1").Value = "Malas comunicaciones"
$ws.Range("D1").Value = "Delincuencia zona"
$ws.Range("E1").Value = "Número viviendas"
$ws.Range("F1").Value = "Contaminación"
$ws.Range("G1").Value = "Falta de servicios de aseo"
$ws.Range("H1").Value = "Aragón"
$ws.Range("I1").Value = "Ruidos exteriores"
